$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.Style = "Normal"
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '27.349.61'
Set-TextValue $ws 'E2' '  -3.21%  '
Set-TextValue $ws 'D3' '1.813.26'
Set-TextValue $ws 'E3' '  -3.63%  '
Set-TextValue $ws 'D4' '1.003'
Set-TextValue $ws 'E4' '  +0.01%  '
Set-TextValue $ws 'D5' '311.17'
Set-TextValue $ws 'E5' '  -1.77%  '
Set-TextValue $ws 'E6' '  +0.01%  '
Set-TextValue $ws 'D7' '0.4217'
Set-TextValue $ws 'E7' '  -2.38%  '
Set-TextValue $ws 'D8' '0.3566'
Set-TextValue $ws 'E8' '  -3.92%  '
Set-TextValue $ws 'D9' '0.07168'
Set-TextValue $ws 'E9' '  -3.45%  '
Set-TextValue $ws 'D10' '0.8491'
Set-TextValue $ws 'E10' '  -4.41%  '
Set-TextValue $ws 'B11' 'WrappedEther'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws 'D11' '1.929.78'
Set-TextValue $ws 'E11' '  +2.62%  '
Set-TextValue $ws 'B12' 'Solana'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws 'D12' '20.24'
Set-TextValue $ws 'E12' '  -4.53%  '
Set-TextValue $ws 'D13' '5.332'
Set-TextValue $ws 'E13' '  -2.87%  '
Set-TextValue $ws 'D14' '6.389'
Set-TextValue $ws 'E14' '  -3.94%  '
Set-TextValue $ws 'D15' '0.06926'
Set-TextValue $ws 'E15' '  -0.79%  '
Set-TextValue $ws 'E16' '  +0.02%  '
Set-TextValue $ws 'D17' '81.58'
Set-TextValue $ws 'E17' '  +0.26%  '
Set-TextValue $ws 'D18' '0.000008840'
Set-TextValue $ws 'E18' '  -3.34%  '
Set-TextValue $ws 'D19' '1.003'
Set-TextValue $ws 'E19' '  +0.06%  '
Set-TextValue $ws 'D20' '15.17'
Set-TextValue $ws 'E20' '  -3.25%  '
Set-TextValue $ws 'D21' '27.559.93'
Set-TextValue $ws 'E21' '  -2.89%  '
Set-TextValue $ws 'D22' '5.099'
Set-TextValue $ws 'E22' '  -0.11%  '
Set-TextValue $ws 'D23' '10.96'
Set-TextValue $ws 'E23' '  -0.71%  '
Set-TextValue $ws 'D24' '2.086.75'
Set-TextValue $ws 'E24' '  -2.29%  '
Set-TextValue $ws 'D25' '1.970'
Set-TextValue $ws 'E25' '  -1.23%  '
Set-TextValue $ws 'D26' '153.88'
Set-TextValue $ws 'E26' '  -0.22%  '
Set-TextValue $ws 'D27' '18.25'
Set-TextValue $ws 'E27' '  -2.93%  '
Set-TextValue $ws 'D28' '5.101'
Set-TextValue $ws 'E28' '  -6.66%  '
Set-TextValue $ws 'D29' '113.41'
Set-TextValue $ws 'E29' '  -4.20%  '
Set-TextValue $ws 'D30' '1.731'
Set-TextValue $ws 'E30' '  -9.20%  '
Set-TextValue $ws 'D31' '0.08898'
Set-TextValue $ws 'E31' '  -1.10%  '
Set-TextValue $ws 'D32' '0.7428'
Set-TextValue $ws 'E32' '  -7.25%  '
Set-TextValue $ws 'D33' '4.487'
Set-TextValue $ws 'E33' '  -4.51%  '
Set-TextValue $ws 'D34' '2.922'
Set-TextValue $ws 'E34' '  -2.57%  '
Set-TextValue $ws 'D35' '1.114'
Set-TextValue $ws 'E35' '  -5.79%  '
Set-TextValue $ws 'D36' '1.002'
Set-TextValue $ws 'E36' '  -0.01%  '
Set-TextValue $ws 'D37' '1.074'
Set-TextValue $ws 'E37' '  -5.55%  '
Set-TextValue $ws 'D38' '0.05210'
Set-TextValue $ws 'E38' '  -5.14%  '
Set-TextValue $ws 'D39' '0.01910'
Set-TextValue $ws 'E39' '  -2.86%  '
Set-TextValue $ws 'D40' '2.771'
Set-TextValue $ws 'E40' '  -4.78%  '
Set-TextValue $ws 'D41' '0.1645'
Set-TextValue $ws 'E41' '  -3.43%  '
Set-TextValue $ws 'D42' '0.5001'
Set-TextValue $ws 'E42' '  -3.59%  '
Set-TextValue $ws 'D43' '6.315'
Set-TextValue $ws 'E43' '  -8.45%  '
Set-TextValue $ws 'D44' '8.236'
Set-TextValue $ws 'E44' '  -4.63%  '
Set-TextValue $ws 'B45' 'EnergySwap'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D45' '10.30'
Set-TextValue $ws 'E45' '  -3.08%  '
Set-TextValue $ws 'B46' 'Quant'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D46' '105.18'
Set-TextValue $ws 'E46' '  -0.83%  '
Set-TextValue $ws 'D47' '0.06420'
Set-TextValue $ws 'E47' '  -3.02%  '
Set-TextValue $ws 'D48' '1.002'
Set-TextValue $ws 'E48' '  -0.03%  '
Set-TextValue $ws 'D49' '0.4608'
Set-TextValue $ws 'E49' '  -3.69%  '
Set-TextValue $ws 'D50' '1.606'
Set-TextValue $ws 'E50' '  -3.63%  '
Set-TextValue $ws 'D51' '63.35'
Set-TextValue $ws 'E51' '  -3.18%  '

Write-Output "Applied cryptos update"
